# Apply the "fixing figs tables for new habglm" edit:
#  1. Refresh the cached "datetimeFigureOut" field text (3/29/2019 -> 4/7/2019)
#     on the slide master and every slide layout's Date placeholder.
#  2. Append " metrics" to the "CRAM" / "IPI" labels on Slide 1's
#     "Rectangle: Rounded Corners 11" shape (inside "Group 1").

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Length -gt 0) {
                $full = $tr.Characters(1, $tr.Length)
                $full.Text = "4/7/2019"
            }
        }
    }
}

# --- Slide master ---
Update-DatePlaceholder $p.SlideMaster.Shapes

# --- Every slide layout ---
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# --- Slide 1: "CRAM" -> "CRAM metrics", "IPI" -> "IPI metrics" ---
$slide1 = $p.Slides.Item(1)
$group1 = $slide1.Shapes.Item("Group 1")
$box = $group1.GroupItems.Item(2)
$tr = $box.TextFrame.TextRange

$cram = $tr.Characters(1, 4)
$cram.Text = "CRAM metrics"

$ipi = $tr.Characters(14, 3)
$ipi.Text = "IPI metrics"
